# Add the "2022-Q1" worksheet (positioned between "2021-Q4" and "总计") and
# populate it with that quarter's fund-holding data, then refresh the
# "总计" (totals) sheet with a new summary row for 2022-Q1.
#
# NOTE: worksheet object references fetched *before* a Worksheets.Add() call
# can end up rebound to the newly inserted sheet in this runtime, so every
# sheet handle used after an Add() is (re)fetched by name right before use.

$wb = $excel.ActiveWorkbook

# --- Create the new "2022-Q1" sheet positioned right before "总计" ---------
$wsTotalBeforeAdd = $wb.Worksheets.Item("总计")
$wsQ1 = $wb.Worksheets.Add($wsTotalBeforeAdd)
$wsQ1.Name = "2022-Q1"

# Re-fetch by name defensively (Add() can rebind stale handles).
$wsQ1 = $wb.Worksheets.Item("2022-Q1")

# Header row (bold, centered, thin border - matches the other sheets' style)
$headerRange = $wsQ1.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$wsQ1.Range("B1").Value = "基金代码"
$wsQ1.Range("C1").Value = "基金名称"
$wsQ1.Range("D1").Value = "基金规模"
$wsQ1.Range("E1").Value = "股票总仓位"
$wsQ1.Range("F1").Value = "仓位占比"
$wsQ1.Range("G1").Value = "持有市值(亿元)"
$wsQ1.Range("H1").Value = "仓位排名"

# Column A (row index) gets the same bold/centered/bordered style as the
# other sheets use for their index column.
$colA = $wsQ1.Range("A2:A6")
$colA.Font.Bold = $true
$colA.HorizontalAlignment = -4108
$colA.VerticalAlignment = -4160
$colA.Borders.LineStyle = 1

# Columns B-G hold text values (fund codes/names/numbers-as-text); force
# text number format so numeric-looking strings (e.g. fund codes with
# leading zeros, "112.21") are not silently coerced into numeric cells.
$wsQ1.Range("B2:G6").NumberFormat = "@"

$data = @(
    @(0, "003567", "华夏行业景气混合",               "112.21", "91.63", "2.76", "3.0970", 7),
    @(1, "001907", "国投瑞银境煊灵活配置混合A",       "2.61",   "90.44", "4.72", "0.1232", 5),
    @(2, "001908", "国投瑞银境煊灵活配置混合C",       "1.75",   "90.44", "4.72", "0.0826", 5),
    @(3, "015309", "国投瑞银境煊灵活配置混合E",       "0.33",   "90.44", "4.72", "0.0156", 5),
    @(4, "005281", "中科沃土转型升级灵活配置混合",     "0.10",   "21.75", "1.29", "0.0013", 7)
)

$r = 2
foreach ($row in $data) {
    $wsQ1.Cells.Item($r, 1).Value = $row[0]
    $wsQ1.Cells.Item($r, 2).Value = $row[1]
    $wsQ1.Cells.Item($r, 3).Value = $row[2]
    $wsQ1.Cells.Item($r, 4).Value = $row[3]
    $wsQ1.Cells.Item($r, 5).Value = $row[4]
    $wsQ1.Cells.Item($r, 6).Value = $row[5]
    $wsQ1.Cells.Item($r, 7).Value = $row[6]
    $wsQ1.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# --- Refresh the "总计" sheet: add a new row for 2022-Q1 above 2021-Q4 -----
$wsTotal = $wb.Worksheets.Item("总计")

$wsTotal.Range("B2:B3").NumberFormat = "@"

# Row 2: new 2022-Q1 summary
$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q1"
$wsTotal.Range("C2").Value = 5
$wsTotal.Range("D2").Value = 3.32

# Row 3: existing 2021-Q4 summary, shifted down one row
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2021-Q4"
$wsTotal.Range("C3").Value = 3
$wsTotal.Range("D3").Value = 4.9

# Column A keeps the bold/centered/bordered style for every data row. A2
# already carried that style before this edit (pre-existing "0" row); only
# the newly-shifted A3 needs it applied explicitly.
$wsTotalA3 = $wsTotal.Range("A3")
$wsTotalA3.Font.Bold = $true
$wsTotalA3.HorizontalAlignment = -4108
$wsTotalA3.VerticalAlignment = -4160
$wsTotalA3.Borders.LineStyle = 1
